# Word COM-interop script implementing the MonsterSlay doc edit:
#  - Rewrite the "One issue..." paragraph to talk about weapon-switching.
#  - Add a new "Some other things are:" paragraph followed by a bulleted
#    list of four items, the last one ending in three manual line breaks.

$d = $word.ActiveDocument

$quote_open  = [char]8220
$quote_close = [char]8221

# --- 1. Rewrite paragraph 7 ("One issue that needs to be addressed ...") ---
$p7 = $d.Paragraphs.Item(7)
$rng = $p7.Range
$rng.MoveEnd(1, -1)
$rng.Text = "One issue that needs to be addressed is the switching weapon-functionality. This seems a bit buggy, as when you press " + $quote_open + "w" + $quote_close + ", the bow or sword will not always be chosen."

# --- 2. Insert the new paragraphs after paragraph 7 -------------------------
# First a blank paragraph right after paragraph 7.
$ip = $d.Paragraphs.Item(7).Range
$ip.Collapse(0)
$ip.InsertParagraphAfter()

# Then five more blank paragraphs, one for each remaining new line of text.
$newTexts = @("Some other things are:", `
              "Adding pictures or animations to the player/monster/power-up", `
              "Create a new weapon that tracks and follows the monsters", `
              "Make changes to the game platform", `
              "Create new maps")

$curIndex = 8
foreach ($t in $newTexts) {
    $p = $d.Paragraphs.Item($curIndex)
    $ip2 = $p.Range
    $ip2.Collapse(0)
    $ip2.InsertParagraphAfter()
    $curIndex = $curIndex + 1
}

# Fill in the text for paragraphs 9..13.
$idx = 9
foreach ($t in $newTexts) {
    $pr = $d.Paragraphs.Item($idx).Range
    $pr.MoveEnd(1, -1)
    $pr.Text = $t
    $idx = $idx + 1
}

# --- 3. Turn paragraphs 10-13 into a shared bulleted list -------------------
# Apply the "List Paragraph" style to each item first ...
for ($i = 10; $i -le 13; $i++) {
    $d.Paragraphs.Item($i).Style = "Listeavsnitt"
}

# ... then apply one shared bullet-list definition across all four at once so
# they end up pointing at the same numId.
$liStart = $d.Paragraphs.Item(10).Range.Start
$liEnd = $d.Paragraphs.Item(13).Range.End
$listRng = $d.Range($liStart, $liEnd)
$listRng.ListFormat.ApplyBulletDefault()

# Bring the auto-generated "List Paragraph" style in line with the Word
# built-in definition (localized id, canonical name, indent, priority).
$liStyle = $d.Styles.Item("Listeavsnitt")
$liStyle.Priority = 34
$liStyle.NameLocal = "List Paragraph"
$liStyle.ParagraphFormat.LeftIndent = 36
$liStyle.NoSpaceBetweenParagraphsOfSameStyle = $true

# --- 4. Add the three trailing manual line breaks on the last bullet --------
$lastItem = $d.Paragraphs.Item(13).Range
$lastItem.Collapse(0)
$lastItem.InsertBreak(6)
$lastItem.Collapse(0)
$lastItem.InsertBreak(6)
$lastItem.Collapse(0)
$lastItem.InsertBreak(6)

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
